# fix: changed data type to text
#
# On the "attributes" sheet (sheet3 / rd3stats.xlsx model), the redundant
# "id" attribute row describing rd3stats_treedata was removed, and the
# idAttribute/nillable flags that used to live on that row were moved onto
# the "subjectID" row (which is the real id attribute for that entity).
# Separately, the "json" attribute's dataType was changed from "string" to
# "text".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Remove the old "id" row (row 2) for rd3stats_treedata; everything below
# shifts up one row.
$ws.Rows(2).Delete()

# The "id" row used to be flagged as the idAttribute (TRUE) and non-nillable
# (FALSE); re-apply those flags to the row that is now in its place
# (subjectID), which is the attribute that truly acts as the id.
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $false

# The "json" attribute row (now row 4 after the shift) should be typed as
# "text" rather than "string".
$ws.Range("D4").Value = "text"

# Reflect the author's final selection on the sheet.
$ws.Activate()
$ws.Range("D5").Select()
